$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New inventory rows for September data (added by infotech), appended after
# the existing 301 data rows.
$newRows = @(
    @{ Row = 302; Cat = "DRINK"; Sub = "CHAMPAGNES";          Item = "DOM PERIGNON" },
    @{ Row = 303; Cat = "DRINK"; Sub = "COGNAC";              Item = "HENNESSY 35CL" },
    @{ Row = 304; Cat = "DRINK"; Sub = "COGNAC";              Item = "REMY MARTIN VSOP BTL" },
    @{ Row = 305; Cat = "DRINK"; Sub = "SHOTS ET DIGESTIFS";  Item = "JET 27" },
    @{ Row = 306; Cat = "DRINK"; Sub = "SHOTS ET DIGESTIFS";  Item = "JET 27 BTL" },
    @{ Row = 307; Cat = "DRINK"; Sub = "SOFTS";               Item = "THE" }
)

foreach ($r in $newRows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.Cat
    $ws.Cells.Item($r.Row, 2).Value = $r.Sub
    $ws.Cells.Item($r.Row, 3).Value = $r.Item
}

# Rows 302-306 carry an extra "vertical centered" style on the middle
# (sub-category) column; row 307 keeps the default style.
$ws.Range("B302:B306").VerticalAlignment = -4108

# Update the active selection to mirror where the user ended up after typing
# the new rows.
$ws.Range("B313").Select() | Out-Null
